$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Tnfsf15"
$ws.Range("C2").Value = "Tnfrsf25"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.4017626666666667
$ws.Range("H2").Value = 1.205288
$ws.Range("I2").Value = 0.2702221345362258
$ws.Range("J2").Value = 0.2702221345362258
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 2.644859666666667
$ws.Range("N2").Value = 7.934579
$ws.Range("O2").Value = 0.250317448389438
$ws.Range("P2").Value = 0.250317448389438
$ws.Range("Q2").Value = 1.062605872639111
$ws.Range("R2").Value = 9.563452853752
$ws.Range("S2").Value = 0.06764131521545547
$ws.Range("T2").Value = 0.06764131521545548

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Tnfsf15"
$ws.Range("C3").Value = "Tnfrsf25"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.4017626666666667
$ws.Range("H3").Value = 1.205288
$ws.Range("I3").Value = 0.2702221345362258
$ws.Range("J3").Value = 0.2702221345362258
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.7341406666666667
$ws.Range("N3").Value = 2.202422
$ws.Range("O3").Value = 0.06948127371556359
$ws.Range("P3").Value = 0.06948127371556359
$ws.Range("Q3").Value = 0.2949503119484445
$ws.Range("R3").Value = 2.654552807536
$ws.Range("S3").Value = 0.01877537809371535
$ws.Range("T3").Value = 0.01877537809371535

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Tnfsf15"
$ws.Range("C4").Value = "Tnfrsf25"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.4017626666666667
$ws.Range("H4").Value = 1.205288
$ws.Range("I4").Value = 0.2702221345362258
$ws.Range("J4").Value = 0.2702221345362258
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1398853333333333
$ws.Range("N4").Value = 0.419656
$ws.Range("O4").Value = 0.0132391673359504
$ws.Range("P4").Value = 0.0132391673359504
$ws.Range("Q4").Value = 0.05620070454755555
$ws.Range("R4").Value = 0.5058063409279999
$ws.Range("S4").Value = 0.003577516057002795
$ws.Range("T4").Value = 0.003577516057002795

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Tnfsf15"
$ws.Range("C5").Value = "Tnfrsf25"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.4017626666666667
$ws.Range("H5").Value = 1.205288
$ws.Range("I5").Value = 0.2702221345362258
$ws.Range("J5").Value = 0.2702221345362258
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 7.047136333333333
$ws.Range("N5").Value = 21.141409
$ws.Range("O5").Value = 0.666962110559048
$ws.Range("P5").Value = 0.6669621105590481
$ws.Range("Q5").Value = 2.831276285643555
$ws.Range("R5").Value = 25.481486570792
$ws.Range("S5").Value = 0.1802279251700522
$ws.Range("T5").Value = 0.1802279251700522

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Tnfsf15"
$ws.Range("C6").Value = "Tnfrsf25"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 1.085024
$ws.Range("H6").Value = 3.255072
$ws.Range("I6").Value = 0.7297778654637742
$ws.Range("J6").Value = 0.7297778654637743
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 2.644859666666667
$ws.Range("N6").Value = 7.934579
$ws.Range("O6").Value = 0.250317448389438
$ws.Range("P6").Value = 0.250317448389438
$ws.Range("Q6").Value = 2.869736214965333
$ws.Range("R6").Value = 25.827625934688
$ws.Range("S6").Value = 0.1826761331739825
$ws.Range("T6").Value = 0.1826761331739826

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Tnfsf15"
$ws.Range("C7").Value = "Tnfrsf25"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 1.085024
$ws.Range("H7").Value = 3.255072
$ws.Range("I7").Value = 0.7297778654637742
$ws.Range("J7").Value = 0.7297778654637743
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.7341406666666667
$ws.Range("N7").Value = 2.202422
$ws.Range("O7").Value = 0.06948127371556359
$ws.Range("P7").Value = 0.06948127371556359
$ws.Range("Q7").Value = 0.7965602427093333
$ws.Range("R7").Value = 7.169042184384
$ws.Range("S7").Value = 0.05070589562184824
$ws.Range("T7").Value = 0.05070589562184825

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Tnfsf15"
$ws.Range("C8").Value = "Tnfrsf25"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 1.085024
$ws.Range("H8").Value = 3.255072
$ws.Range("I8").Value = 0.7297778654637742
$ws.Range("J8").Value = 0.7297778654637743
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.1398853333333333
$ws.Range("N8").Value = 0.419656
$ws.Range("O8").Value = 0.0132391673359504
$ws.Range("P8").Value = 0.0132391673359504
$ws.Range("Q8").Value = 0.1517789439146667
$ws.Range("R8").Value = 1.366010495232
$ws.Range("S8").Value = 0.009661651278947605
$ws.Range("T8").Value = 0.009661651278947608

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Tnfsf15"
$ws.Range("C9").Value = "Tnfrsf25"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 1.085024
$ws.Range("H9").Value = 3.255072
$ws.Range("I9").Value = 0.7297778654637742
$ws.Range("J9").Value = 0.7297778654637743
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 7.047136333333333
$ws.Range("N9").Value = 21.141409
$ws.Range("O9").Value = 0.666962110559048
$ws.Range("P9").Value = 0.6669621105590481
$ws.Range("Q9").Value = 7.646312052938667
$ws.Range("R9").Value = 68.816808476448
$ws.Range("S9").Value = 0.4867341853889958
$ws.Range("T9").Value = 0.486734185388996

